$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the summary block above the table
# ---------------------------------------------------------------------------
# VALOR MORA total: 227760 -> 284700 (one more period added, 56940 * 5)
$ws.Range("E11").Value = 284700
# Cant. Periodos: 4 -> 5
$ws.Range("F13").Value = 5

# ---------------------------------------------------------------------------
# 2. Re-sort the existing "Periodo Mora" rows (16-18) into ascending order.
#    Before: row16=2507, row17=2506, row18=2505, row19=2504
#    After : row16=2504, row17=2505, row18=2506, row19=2507 (row19 is new)
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"

# ---------------------------------------------------------------------------
# 3. Insert a new row before the current last data row (row 19), shifting the
#    old row 19 (period 2504) down to row 20, and everything below it as well
#    (the signature block moves from rows 24-25 to rows 25-26).
# ---------------------------------------------------------------------------
$ws.Rows("19:19").Insert()

# Copy formatting only (borders/fill/font/number format) from the row above
# (row 18) onto the freshly inserted row 19, matching the original table's
# "middle row" style.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row 19 with period 2507 data (same worker/values pattern).
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "79940828"
$ws.Range("D19").Value = "CARLOS ANDRES CABALLERO PULGARIN"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20 (old row 19, shifted down) already has the correct style and the
# worker data; it now represents the newly added period 2508.
$ws.Range("E20").Value = "2508"

